$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.759.58'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.887.03'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7927'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.41'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3168'
$ws.Range('E8').Value = '  +1.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.49'
$ws.Range('E9').Value = '  -2.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07035'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08046'
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7675'
$ws.Range('E12').Value = '  +3.99%  '
$ws.Range('D13').Value = '1.879.62'
$ws.Range('E13').Value = '  -1.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.293'
$ws.Range('E14').Value = '  +2.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.92'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = '29.774.54'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.78'
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.926'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.32'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007703'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.157.59'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.101'
$ws.Range('E23').Value = '  +17.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1625'
$ws.Range('E25').Value = '  +11.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.288'
$ws.Range('E26').Value = '  +1.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.85'
$ws.Range('E27').Value = '  -2.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.62'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.057'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.370'
$ws.Range('E30').Value = '  +1.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.534'
$ws.Range('E31').Value = '  +1.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.437'
$ws.Range('E32').Value = '  +3.86%  '
$ws.Range('E33').Value = '  +2.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.082'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.261'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7340'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.717'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01918'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.767'
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4406'
$ws.Range('E41').Value = '  +0.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '71.93'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8399'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').Value = '1.030.84'
$ws.Range('E46').Value = '  +5.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.93'
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.856'
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.875'
$ws.Range('E49').Value = '  +1.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.433'
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').Value = '2.040.97'
$ws.Range('E51').Value = '  -0.75%  '
